$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-26 Thursday", 2)

$d.Content.Find.Execute("132÷7=18, 6", $true, $false, $false, $false, $false, $true, 1, $false, "697÷7=99, 4", 2)
$d.Content.Find.Execute("557÷3=185, 2", $true, $false, $false, $false, $false, $true, 1, $false, "543÷5=108, 3", 2)
$d.Content.Find.Execute("317÷5=63, 2", $true, $false, $false, $false, $false, $true, 1, $false, "892÷7=127, 3", 2)
$d.Content.Find.Execute("676÷3=225, 1", $true, $false, $false, $false, $false, $true, 1, $false, "999÷9=111, 0", 2)
$d.Content.Find.Execute("918÷8=114, 6", $true, $false, $false, $false, $false, $true, 1, $false, "450÷7=64, 2", 2)

$d.Content.Find.Execute("287÷2=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "915÷9=101, 6", 2)
$d.Content.Find.Execute("244÷8=30, 4", $true, $false, $false, $false, $false, $true, 1, $false, "854÷9=94, 8", 2)
$d.Content.Find.Execute("146÷7=20, 6", $true, $false, $false, $false, $false, $true, 1, $false, "330÷4=82, 2", 2)
$d.Content.Find.Execute("972÷7=138, 6", $true, $false, $false, $false, $false, $true, 1, $false, "925÷3=308, 1", 2)
$d.Content.Find.Execute("639÷9=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "998÷8=124, 6", 2)

$d.Content.Find.Execute("281÷7=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "782÷8=97, 6", 2)
$d.Content.Find.Execute("815÷4=203, 3", $true, $false, $false, $false, $false, $true, 1, $false, "835÷4=208, 3", 2)
$d.Content.Find.Execute("275÷3=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "448÷8=56, 0", 2)
$d.Content.Find.Execute("413÷7=59, 0", $true, $false, $false, $false, $false, $true, 1, $false, "794÷9=88, 2", 2)
$d.Content.Find.Execute("429÷9=47, 6", $true, $false, $false, $false, $false, $true, 1, $false, "206÷8=25, 6", 2)

$d.Content.Find.Execute("730÷5=146, 0", $true, $false, $false, $false, $false, $true, 1, $false, "833÷7=119, 0", 2)
$d.Content.Find.Execute("411÷3=137, 0", $true, $false, $false, $false, $false, $true, 1, $false, "298÷6=49, 4", 2)
$d.Content.Find.Execute("889÷5=177, 4", $true, $false, $false, $false, $false, $true, 1, $false, "363÷4=90, 3", 2)
$d.Content.Find.Execute("830÷2=415, 0", $true, $false, $false, $false, $false, $true, 1, $false, "131÷8=16, 3", 2)
$d.Content.Find.Execute("618÷8=77, 2", $true, $false, $false, $false, $false, $true, 1, $false, "820÷3=273, 1", 2)

$d.Content.Find.Execute("943÷3=314, 1", $true, $false, $false, $false, $false, $true, 1, $false, "137÷8=17, 1", 2)
$d.Content.Find.Execute("742÷2=371, 0", $true, $false, $false, $false, $false, $true, 1, $false, "823÷6=137, 1", 2)
$d.Content.Find.Execute("231÷5=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "271÷7=38, 5", 2)
$d.Content.Find.Execute("734÷8=91, 6", $true, $false, $false, $false, $false, $true, 1, $false, "369÷5=73, 4", 2)
$d.Content.Find.Execute("173÷9=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "940÷8=117, 4", 2)
